$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 704, shifting existing rows 704:806 down to 705:807
$ws.Rows.Item(704).Insert()

# Populate the newly inserted row 704 with the new data entry
$ws.Range("A704").Value = 10
$ws.Range("B704").Value = "Vega Modelo de Temuco"
$ws.Range("C704").Value = "La Araucanía"
$ws.Range("D704").Value = 45034
$ws.Range("E704").Value = 9
$ws.Range("F704").Value = 100112045
$ws.Range("G704").Value = "Zapallo"
$ws.Range("H704").Value = "Camote"
$ws.Range("I704").Value = "1a (cosecha)"
$ws.Range("J704").Value = 580
$ws.Range("K704").Value = 600
$ws.Range("L704").Value = 600
$ws.Range("M704").Value = 600
$ws.Range("N704").Value = "$/kilo (volumen en unidades)"
$ws.Range("O704").Value = "Región del Maule"
$ws.Range("P704").Value = 600
$ws.Range("Q704").Value = 1
$ws.Range("R704").Value = "Hortaliza"
